$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: theta_se values (previously all "(nan)")
$ws.Range("B4").Value = "(0.0)"
$ws.Range("C4").Value = "(0.1)"
$ws.Range("D4").Value = "(0.13)"
$ws.Range("E4").Value = "(0.03)"
$ws.Range("F4").Value = "(0.39)"
$ws.Range("G4").Value = "(0.24)"

# Row 6: lambda_se values (previously all "(nan)")
$ws.Range("B6").Value = "(0.0)"
$ws.Range("C6").Value = "(0.01)"
$ws.Range("D6").Value = "(0.01)"
$ws.Range("E6").Value = "(0.06)"
$ws.Range("F6").Value = "(0.21)"
$ws.Range("G6").Value = "(0.04)"
